$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text interpretation so that numeric-
# looking strings (prices such as "5.52" or padded percentages) are stored
# verbatim as text instead of being parsed into floating point numbers.
# The temporary "@" (Text) number format is reverted to the Normal style
# right after the write so no stray cell formatting is left behind.
function Set-TextValue($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "65.023.82"
Set-TextValue "E2" "  +1.38%  "
Set-TextValue "D3" "3.181.89"
Set-TextValue "E3" "  +1.65%  "
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "617.24"
Set-TextValue "E5" "  +1.89%  "
Set-TextValue "D6" "148.07"
Set-TextValue "E6" "  -0.34%  "
Set-TextValue "E7" "  -0.08%  "
Set-TextValue "D8" "3.177.17"
Set-TextValue "E8" "  +1.54%  "
Set-TextValue "E9" "  +0.06%  "
Set-TextValue "E10" "  +0.86%  "
Set-TextValue "D11" "5.52"
Set-TextValue "E11" "  -1.19%  "
Set-TextValue "E12" "  +0.02%  "
Set-TextValue "D13" "0.0000263"
Set-TextValue "E13" "  +1.97%  "
Set-TextValue "D14" "36.06"
Set-TextValue "E14" "  -1.58%  "
Set-TextValue "D15" "3.699.37"
Set-TextValue "E15" "  +1.42%  "
Set-TextValue "E16" "  +3.35%  "
Set-TextValue "D17" "64.998.70"
Set-TextValue "E17" "  +1.23%  "
Set-TextValue "D18" "3.176.92"
Set-TextValue "E18" "  +1.25%  "
Set-TextValue "E19" "  -0.33%  "
Set-TextValue "D20" "482.04"
Set-TextValue "E20" "  +0.58%  "
Set-TextValue "D21" "14.79"
Set-TextValue "E21" "  +1.59%  "
Set-TextValue "D22" "0.724"
Set-TextValue "E22" "  +2.21%  "
Set-TextValue "E23" "  +3.38%  "
Set-TextValue "E24" "  +1.39%  "
Set-TextValue "D25" "84.77"
Set-TextValue "E25" "  +1.20%  "
Set-TextValue "E26" "  +0.02%  "
Set-TextValue "B27" "RenderToken"
Set-TextValue "C27" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D27" "8.74"
Set-TextValue "E27" "  +2.62%  "
Set-TextValue "B28" "PancakeSwap"
Set-TextValue "C28" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D28" "2.84"
Set-TextValue "E28" "  -3.00%  "
Set-TextValue "B29" "NEARProtocol"
Set-TextValue "C29" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D29" "7.07"
Set-TextValue "E29" "  +1.70%  "
Set-TextValue "B30" "Hedera"
Set-TextValue "C30" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D30" "0.120"
Set-TextValue "E30" "  -4.91%  "
Set-TextValue "B31" "ImmutableX"
Set-TextValue "C31" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D31" "2.12"
Set-TextValue "E31" "  -4.94%  "
Set-TextValue "E32" "  +0.04%  "
Set-TextValue "E33" "  -0.30%  "
Set-TextValue "D34" "26.84"
Set-TextValue "E34" "  +0.84%  "
Set-TextValue "E35" "  +2.90%  "
Set-TextValue "D36" "0.0₃0798"
Set-TextValue "E36" "  +7.16%  "
Set-TextValue "D37" "6.07"
Set-TextValue "E37" "  -0.09%  "
Set-TextValue "D38" "3.21"
Set-TextValue "E38" "  -0.37%  "
Set-TextValue "D39" "53.18"
Set-TextValue "E39" "  -2.47%  "
Set-TextValue "D40" "470.36"
Set-TextValue "E40" "  +4.85%  "
Set-TextValue "D41" "0.0402"
Set-TextValue "E41" "  +0.84%  "
Set-TextValue "E42" "  -2.26%  "
Set-TextValue "E43" "  -0.42%  "
Set-TextValue "D44" "2.866.81"
Set-TextValue "E44" "  -0.16%  "
Set-TextValue "D45" "2.36"
Set-TextValue "E45" "  +2.25%  "
Set-TextValue "D46" "0.271"
Set-TextValue "E46" "  +0.07%  "
Set-TextValue "D47" "2.47"
Set-TextValue "E47" "  +7.21%  "
Set-TextValue "B48" "Arweave"
Set-TextValue "C48" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D48" "37.74"
Set-TextValue "E48" "  +13.69%  "
Set-TextValue "B49" "InjectiveProtocol"
Set-TextValue "C49" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D49" "26.93"
Set-TextValue "E49" "  +1.51%  "
Set-TextValue "D50" "0.999"
Set-TextValue "E50" "  +0.06%  "
Set-TextValue "E51" "  +0.00%  "
